# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value()
$newText = $oldText.Replace(
    "1000 Bs = 6.95 = 27198.83 pesos",
    "1000 Bs = 6.99 = 27378.58 pesos"
)
$newText = $newText.Replace(
    "27198.83 pesos = 6.91 = 966.95 Bs",
    "27378.58 pesos = 6.97 = 974.68 Bs"
)
$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 142.995
$wsTasas.Range("N12").Value = 3929.75
$wsTasas.Range("O12").Value = 139.9
